# Regenerate the localization-status report:
#  - the "zh-cn" / "de-de" status columns on the Overview sheet, and the
#    "Status" column on each per-language sheet, move from
#    "Ready for handoff" to "In Translation"
#  - the now-shorter status text means those columns get narrower once the
#    report is (re)generated / autofit

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # compare with the literal on the left so PowerShell doesn't coerce
        # the literal to the cell's (possibly boolean) type
        if ($oldStatus -eq $cell.Value2) {
            $cell.Value = $newStatus
        }
    }
}

# Overview sheet: zh-cn (E) and de-de (F) status columns narrow to fit the
# new, shorter status text
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn / de-de sheets: Status column (C) narrows the same way
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
